$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task descriptions for weeks 2-10 (rows 3,4,6,7,9,10,12,13,15) ---
# Row 3 (week 2) - wraps, row height 30
$ws.Cells.Item(3, 2).Value = "Page de création de milieu de stage"
$ws.Cells.Item(3, 3).Value = "Page de création de stage"
$ws.Cells.Item(3, 4).Value = "Page de liste des milieu de stages et de stage"

# Row 4 (week 3) - wraps, row height 30
$ws.Cells.Item(4, 2).Value = "Page de modification et duplication de milieu de stage"
$ws.Cells.Item(4, 3).Value = "Page de modification et duplication de stage"
$ws.Cells.Item(4, 4).Value = "Impression de liste de stage + documentation scolaire"

# Row 6 (week 4) - "Retour sprint 1" in all 3 columns
$ws.Cells.Item(6, 2).Value = "Retour sprint 1"
$ws.Cells.Item(6, 3).Value = "Retour sprint 1"
$ws.Cells.Item(6, 4).Value = "Retour sprint 1"

# Row 7 (week 5) - wraps, row height 30
$ws.Cells.Item(7, 2).Value = "Ordre de préférance des stages de l'étudiant"
$ws.Cells.Item(7, 3).Value = "Documentation scolaire + liste des stage vue par l'étudiant"
$ws.Cells.Item(7, 4).Value = "Filtres de la liste des stages présentés à l'étudiant"

# Row 9 (week 6) - "Retour sprint 2" in all 3 columns
$ws.Cells.Item(9, 2).Value = "Retour sprint 2"
$ws.Cells.Item(9, 3).Value = "Retour sprint 2"
$ws.Cells.Item(9, 4).Value = "Retour sprint 2"

# Row 10 (week 7) - wraps, row height 60
$ws.Cells.Item(10, 2).Value = "Consulter les détails des étudiants et attribuer un stage en dehors des choix de l'étudiant"
$ws.Cells.Item(10, 3).Value = "Attribuer un superviseur au stage de l'étudiant et visualiser les choix tentatif de l'étudiant"
$ws.Cells.Item(10, 4).Value = "Approuver le choix de stage de l'étudiant"

# Row 12 (week 8) - "Retour sprint 3" in all 3 columns
$ws.Cells.Item(12, 2).Value = "Retour sprint 3"
$ws.Cells.Item(12, 3).Value = "Retour sprint 3"
$ws.Cells.Item(12, 4).Value = "Retour sprint 3"

# Row 13 (week 9) - "Réglage des derniers bugs" in all 3 columns
$ws.Cells.Item(13, 2).Value = "Réglage des derniers bugs"
$ws.Cells.Item(13, 3).Value = "Réglage des derniers bugs"
$ws.Cells.Item(13, 4).Value = "Réglage des derniers bugs"

# Row 15 (week 10) - "Retour sprint 4" in all 3 columns
$ws.Cells.Item(15, 2).Value = "Retour sprint 4"
$ws.Cells.Item(15, 3).Value = "Retour sprint 4"
$ws.Cells.Item(15, 4).Value = "Retour sprint 4"

# --- Style the new content cells like row 2 (wrap text) ---
$ws.Range("B3:D4").WrapText = $true
$ws.Range("B6:D6").WrapText = $true
$ws.Range("B7:D7").WrapText = $true
$ws.Range("B9:D9").WrapText = $true
$ws.Range("B10:D10").WrapText = $true
$ws.Range("B12:D12").WrapText = $true
$ws.Range("B13:D13").WrapText = $true
$ws.Range("B15:D15").WrapText = $true

# --- Turn on word-wrap for the previously-blank separator rows (Neutre style) ---
$ws.Range("B5:D5").WrapText = $true
$ws.Range("B8:D8").WrapText = $true
$ws.Range("B11:D11").WrapText = $true
$ws.Range("B14:D14").WrapText = $true

# --- Row heights for multi-line rows ---
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).RowHeight = 30
$ws.Rows(7).RowHeight = 30
$ws.Rows(10).RowHeight = 60

# --- Selection moved to G10 in the saved file ---
$ws.Range("G10").Select() | Out-Null
